$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wrap the three hint strings (rows 2-4, column B) with the
# "<color=#00CC00>(...)</color>" marker used for in-game hint text.
$ws.Range("B2").Value = " <color=#00CC00>(Some marks on the body can reveal information about the killer.)</color>"
$ws.Range("B3").Value = " <color=#00CC00>(Think back to the two conversations with Doctor Ran.)</color>"
$ws.Range("B4").Value = " <color=#00CC00>(There was one particular clue that points directly to the answer.)</color>"

# Row 3 now wraps onto a second line like rows 2 and 4, so bump its height
# to match (34 points, same as the other two-line rows).
$ws.Rows.Item(3).RowHeight = 34

# Move the active selection from B6 to B12.
$ws.Range("B12").Select()
